$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Family Member Data")

# Scratch cell used to force text-typed (shared-string) values for
# cells whose content looks numeric (e.g. "1", "2", "962359") without
# leaving stray formatting behind on the sheet.
$scratch = $ws.Cells.Item(100, 26)
$scratch.NumberFormat = "@"

function Set-TextValue($cell, [string]$text) {
    $scratch.Value = $text
    $scratch.Copy()
    $cell.PasteSpecial(-4163)
}

# ---- Header row ----
$headers = @("MemberID","FatherID","MotherID","MemberName","NickName","BirthOrder","Origin","NationalityID","ReligionID","Dob","LunarDob","BirthPlace","IsDead","Dod","LunarDod","PlaceOfDeath","GraveSite","Note","Generation","BloodType","Male","CodeID","Image")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# ---- Row 2 ----
$ws.Cells.Item(2,1).Value = 853
$ws.Cells.Item(2,2).Value = 836
$ws.Cells.Item(2,3).Value = 889
Set-TextValue $ws.Cells.Item(2,4) "2"
Set-TextValue $ws.Cells.Item(2,5) "2"
$ws.Cells.Item(2,6).Value = 1
Set-TextValue $ws.Cells.Item(2,7) "1"
$ws.Cells.Item(2,8).Value = 1
$ws.Cells.Item(2,9).Value = 1
$ws.Cells.Item(2,13).Value = 0
$ws.Cells.Item(2,19).Value = 4
$ws.Cells.Item(2,21).Value = 1
Set-TextValue $ws.Cells.Item(2,22) "962359"

# ---- Row 3 ----
$ws.Cells.Item(3,1).Value = 862
$ws.Cells.Item(3,2).Value = 853
Set-TextValue $ws.Cells.Item(3,4) "ff"
Set-TextValue $ws.Cells.Item(3,5) "ff"
$ws.Cells.Item(3,6).Value = 1
Set-TextValue $ws.Cells.Item(3,7) "1"
$ws.Cells.Item(3,8).Value = 1
$ws.Cells.Item(3,9).Value = 1
$ws.Cells.Item(3,13).Value = 0
$ws.Cells.Item(3,19).Value = 5
$ws.Cells.Item(3,21).Value = 1
Set-TextValue $ws.Cells.Item(3,22) "962359"

# ---- Row 4 ----
$ws.Cells.Item(4,1).Value = 873
$ws.Cells.Item(4,2).Value = 862
Set-TextValue $ws.Cells.Item(4,4) "f"
Set-TextValue $ws.Cells.Item(4,5) "f"
$ws.Cells.Item(4,6).Value = 1
Set-TextValue $ws.Cells.Item(4,7) "1"
$ws.Cells.Item(4,8).Value = 1
$ws.Cells.Item(4,9).Value = 1
$ws.Cells.Item(4,13).Value = 0
$ws.Cells.Item(4,19).Value = 6
$ws.Cells.Item(4,21).Value = 1
Set-TextValue $ws.Cells.Item(4,22) "962359"

# ---- Row 5 ----
$ws.Cells.Item(5,1).Value = 889
Set-TextValue $ws.Cells.Item(5,4) "hh"
Set-TextValue $ws.Cells.Item(5,5) "hh"
$ws.Cells.Item(5,6).Value = 1
Set-TextValue $ws.Cells.Item(5,7) "1"
$ws.Cells.Item(5,8).Value = 1
$ws.Cells.Item(5,9).Value = 1
$ws.Cells.Item(5,13).Value = 0
$ws.Cells.Item(5,19).Value = 3
$ws.Cells.Item(5,21).Value = 0
Set-TextValue $ws.Cells.Item(5,22) "962359"

# Remove the scratch helper cell and its formatting.
$scratch.Clear()

$wb.Save()
